# "Insert im BookingRepositoryDB erstellt"
#
# Adds a new task row ("Insert im BookingRepositoryDB erstellen", done,
# 2019-04-05) to the Tasks sheet, right above the final "9." user-story
# summary row, and restores the selection/scroll state on all three sheets
# to match where the author's cursor ended up after the edit.

$wb = $excel.ActiveWorkbook

# --- UserStories sheet: just a cursor move, no data change ---
$wsUserStories = $wb.Worksheets.Item("UserStories")
$wsUserStories.Range("B13").Select()

# --- Priorisierung sheet: just a cursor move / scroll, no data change ---
$wsPrio = $wb.Worksheets.Item("Priorisierung")
$wsPrio.Select()
$wsPrio.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 3
$wsPrio.Range("B17").Select()

# --- Tasks sheet: the actual new-row insertion ---
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsTasks.Select()

# Insert a new row 47 (pushes the old "9." summary row from 48 down to 49),
# copying formatting from the row above as Excel normally does.
$wsTasks.Rows("47").Insert()

$wsTasks.Range("B47").Value = "Insert im BookingRepositoryDB erstellen"
$wsTasks.Range("C47").Value = "done"
$wsTasks.Range("D47").Value = 43560

# Match the final scroll position / selection from the authored edit.
$wsTasks.Range("A37").Select()
$excel.ActiveWindow.ScrollRow = 37
$wsTasks.Range("E47").Select()
